$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 on the "Rules" sheet currently holds the text "R40" (rule name).
# It needs to become the text "1", stored as a shared string (not a number),
# keeping the cell's existing style/format untouched.
#
# A direct `$ws.Range("B11").Value = "1"` would store "1" as a NUMBER (Excel
# auto-detects numeric-looking literals), and forcing text via NumberFormat
# would also allocate a brand-new cell style for B11. Instead, build the text
# "1" in a scratch cell via a formula (TEXT() always returns a string), then
# copy/paste just the VALUE into B11 so its existing style is preserved.
$scratch = $ws.Range("ZZ1")
$scratch.Formula = '=TEXT(1,"0")'
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# Remove the scratch helper column entirely so it leaves no trace in the
# sheet (no leftover cell, dimension stays the same).
$scratch.EntireColumn.Delete()
